$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.664.97"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "3.786.08"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'595.68"
$ws.Range("D6").Value = "'167.00"
$ws.Range("E6").Value = "  +0.51%  "
$ws.Range("D7").Value = "3.770.69"
$ws.Range("E7").Value = "  +0.77%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("E10").Value = "  -0.08%  "
$ws.Range("D11").Value = "'6.30"
$ws.Range("E11").Value = "  -2.11%  "
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("E13").Value = "  -2.83%  "
$ws.Range("D14").Value = "'35.95"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "4.419.80"
$ws.Range("E15").Value = "  +1.08%  "
$ws.Range("D16").Value = "3.793.98"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").Value = "'18.61"
$ws.Range("E17").Value = "  +4.31%  "
$ws.Range("D18").Value = "67.608.96"
$ws.Range("E18").Value = "  -1.03%  "
$ws.Range("E19").Value = "  +0.85%  "
$ws.Range("E20").Value = "  -0.17%  "
$ws.Range("E21").Value = "  -5.74%  "
$ws.Range("D22").Value = "'459.60"
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("E23").Value = "  +0.18%  "
$ws.Range("E24").Value = "  +4.22%  "
$ws.Range("D25").Value = "'83.33"
$ws.Range("E25").Value = "  -0.70%  "
$ws.Range("D26").Value = "'11.99"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").Value = "'2.11"
$ws.Range("E27").Value = "  -2.88%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D30").Value = "3.930.70"
$ws.Range("E30").Value = "  +1.00%  "
$ws.Range("E31").Value = "  +0.48%  "
$ws.Range("E32").Value = "  +3.67%  "
$ws.Range("E33").Value = "  -1.22%  "
$ws.Range("D34").Value = "'29.64"
$ws.Range("E34").Value = "  -0.53%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  -0.76%  "
$ws.Range("D37").Value = "'0.100"
$ws.Range("E37").Value = "  -0.79%  "
$ws.Range("D38").Value = "'3.36"
$ws.Range("E38").Value = "  -1.77%  "
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("D40").Value = "'0.997"
$ws.Range("E40").Value = "  +0.28%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'1.00"
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D44").Value = "'45.52"
$ws.Range("E44").Value = "  +5.87%  "
$ws.Range("D45").Value = "'48.11"
$ws.Range("E45").Value = "  +3.01%  "
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "'150.09"
$ws.Range("E47").Value = "  +4.11%  "
$ws.Range("D48").Value = "'8.32"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").Value = "'393.19"
$ws.Range("E49").Value = "  +1.13%  "
$ws.Range("D50").Value = "'26.68"
$ws.Range("E50").Value = "  +6.64%  "
$ws.Range("D51").Value = "'1.82"
$ws.Range("E51").Value = "  -5.15%  "
